# Project Log.xlsx - add a new log entry row (2019-06-12, "Voice Activity
# Detection research") above the existing "Total (hrs)" summary rows, and
# extend the running-total SUM() ranges to cover it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook's old row 25 ("Total (hrs)" header row) and row 26 (the SUM
# totals row) need to shift down to rows 26/27 to make room for the new
# data row. Inserting a row above the old row 25 does exactly that and
# also carries down the formatting (cell styles) of the row above it,
# which happens to already match what the new data row needs (same style
# ids as every other "Development" row, e.g. row 24).
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row 25 with the new log entry.
$ws.Cells.Item(25, 2).Value = "Development"
$ws.Cells.Item(25, 3).Value = 43628
$ws.Cells.Item(25, 4).Value = 6
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(25, 6).Value = "1) Voice Activity Detection research (https://www.ncbi.nlm.nih.gov/pmc/articles/PMC4142156/)"

# Remarks column wraps text; this entry wraps onto two lines same as the
# other short remarks rows (e.g. row 19/21/22), so the row is taller than
# the default.
$ws.Rows.Item(25).RowHeight = 29

# Old row 26 is now row 27 ("Total (hrs)" SUM formulas). Extend its ranges
# so the new data row (25) is included in the totals.
$ws.Cells.Item(27, 3).Formula = "=SUM(D27:E27)"
$ws.Cells.Item(27, 4).Formula = "=SUM(D3:D25)"
$ws.Cells.Item(27, 5).Formula = "=SUM(E3:E25)"

# Scroll the view down a couple of rows and move the selection to the new
# last remarks cell, matching where the author's cursor ended up.
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$ws.Range("F26").Select()
